# Horarios actualizados Linea 141 - 945
# Update "last scraped" timestamp from 04:15:01 -> 04:34:13 across all
# sheets, refresh the Minutos/rows data on the LP1912 sheet (including 4
# new departures), and refresh the two summary rows on LP1912-215.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:34:13"
$ws1.Range("A3").Value = "Total filas: 14"

$data1 = @(
    @("04:34:13", "04:47", "215_EL PELIGRO",  13, "LP1912"),
    @("04:34:13", "04:53", "11_ETCHEVERRY",   19, "LP1912"),
    @("04:34:13", "05:11", "17_ROMERO",       37, "LP1912"),
    @("04:34:13", "05:22", "23_HERNANDEZ",    48, "LP1912"),
    @("04:34:13", "05:32", "81_EL PELIGRO",   58, "LP1912"),
    @("04:34:13", "05:44", "14_ABASTO",       70, "LP1912"),
    @("04:34:13", "05:52", "17_ROMERO",       78, "LP1912"),
    @("04:34:13", "06:01", "16_SANTA ANA",    87, "LP1912"),
    @("04:34:13", "06:04", "10_OLMOS",        90, "LP1912"),
    @("04:34:13", "06:11", "215A_EL PATO",    97, "LP1912"),
    @("04:34:13", "06:24", "11_ETCHEVERRY",  110, "LP1912"),
    @("04:34:13", "06:27", "23_HERNANDEZ",   113, "LP1912"),
    @("04:34:13", "06:31", "17X38_ROMERO",   117, "LP1912"),
    @("04:34:13", "06:31", "16_SANTA ANA",   117, "LP1912")
)

$row = 6
foreach ($rec in $data1) {
    $ws1.Cells.Item($row, 1).Value = $rec[0]
    $ws1.Cells.Item($row, 2).Value = $rec[1]
    $ws1.Cells.Item($row, 3).Value = $rec[2]
    $ws1.Cells.Item($row, 4).Value = $rec[3]
    $ws1.Cells.Item($row, 5).Value = $rec[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:34:13"

$data2 = @(
    @("04:34:13", "04:47", "215_EL PELIGRO", 13, "LP1912"),
    @("04:34:13", "06:11", "215A_EL PATO",   97, "LP1912")
)

$row = 6
foreach ($rec in $data2) {
    $ws2.Cells.Item($row, 1).Value = $rec[0]
    $ws2.Cells.Item($row, 2).Value = $rec[1]
    $ws2.Cells.Item($row, 3).Value = $rec[2]
    $ws2.Cells.Item($row, 4).Value = $rec[3]
    $ws2.Cells.Item($row, 5).Value = $rec[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:34:13"
